$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$dateFmt = "[$-409]d/mmm/yyyy;@"
$currFmt = '"₹"#,##0;"₹"\-#,##0'

# --- Row 306: add AMOUNT (F306) that was previously missing ---
$ws.Cells.Item(306, 6).Value = 6000
$ws.Cells.Item(306, 6).NumberFormat = $currFmt

# --- Row 310: new entry ---
$ws.Cells.Item(310, 1).Value = 44797
$ws.Cells.Item(310, 1).NumberFormat = $dateFmt
$ws.Cells.Item(310, 2).Value = "KA03MS2872"
$ws.Cells.Item(310, 3).Value = "BEAT"
$ws.Cells.Item(310, 4).Value = "PMS                                      WW"
$ws.Cells.Item(310, 5).Value = "WORK IN PROGRESS"

# --- Row 311: new entry ---
$ws.Cells.Item(311, 1).Value = 44797
$ws.Cells.Item(311, 1).NumberFormat = $dateFmt
$ws.Cells.Item(311, 2).Value = "KA06Z3574"
$ws.Cells.Item(311, 3).Value = "I20"
$ws.Cells.Item(311, 4).Value = "PMS"
$ws.Cells.Item(311, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(311, 6).Value = 3627
$ws.Cells.Item(311, 6).NumberFormat = $currFmt
$ws.Cells.Item(311, 7).Value = "G PAY"

# --- Row 312: new entry ---
$ws.Cells.Item(312, 1).Value = 44797
$ws.Cells.Item(312, 1).NumberFormat = $dateFmt
$ws.Cells.Item(312, 2).Value = "KA50M1697"
$ws.Cells.Item(312, 3).Value = "VERNA"
$ws.Cells.Item(312, 4).Value = "VCMC"
$ws.Cells.Item(312, 5).Value = "WORK IN PROGRESS"

# --- Row 313: new entry ---
$ws.Cells.Item(313, 1).Value = 44797
$ws.Cells.Item(313, 1).NumberFormat = $dateFmt
$ws.Cells.Item(313, 2).Value = "KA04MM4818"
$ws.Cells.Item(313, 3).Value = "RITZ"
$ws.Cells.Item(313, 4).Value = "GENERAL CHECKUP         WW"
$ws.Cells.Item(313, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(313, 6).Value = 1945
$ws.Cells.Item(313, 6).NumberFormat = $currFmt
$ws.Cells.Item(313, 7).Value = "CREDIT"

$ws.Range("G313").Select()
